# Remove the trailing "Ver no Jupiter ..." / copyright footer block that
# used to follow the "LOM3011: Ensaios Mecânicos (Requisito)" paragraph,
# along with the blank paragraph that separated them from it.
#
# Resulting structure right after the requirement line should simply be:
#   LOM3011: Ensaios Mecânicos (Requisito)
#   <blank paragraph>
#   <page-break paragraph>
# (the blank + page-break paragraph pair that already closed the document).

$d = $word.ActiveDocument

# Locate the "LOM3011: Ensaios Mecânicos (Requisito)" paragraph. Matching
# on a plain ASCII prefix avoids any accented-character encoding pitfalls.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "LOM3011*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'LOM3011: Ensaios Mecanicos (Requisito)' paragraph"
}

# The three paragraphs immediately following it are the ones being removed:
#   1) an empty paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) "© 2020 . Contact: ... Creative Commons Attribution"
$firstToRemove = $targetIndex + 1
$lastToRemove = $targetIndex + 3

$start = $d.Paragraphs.Item($firstToRemove).Range.Start
$end = $d.Paragraphs.Item($lastToRemove).Range.End

$d.Range($start, $end).Delete()

Write-Output ("Removed paragraphs " + $firstToRemove + ".." + $lastToRemove + "; ParagraphCount now " + $d.Paragraphs.Count)
